# Bug 5 fixed + endday texture loaded wrong
#
# The bug-tracker table gains a new row (Bug 5: "When the shop is selected
# in the StartOfDayState, it's still being selected in GameplayState"),
# and the stray "_GoBack" bookmark that Word had left in Bug 2's status
# cell moves down to sit in the new Bug 5's status cell instead.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- 1. Strip the "_GoBack" bookmark out of Bug 2's "Fixed" status cell ---
# (table row 3 = header + Bug1 + Bug2)
$plainFixedXml = "<w:p $wNs><w:pPr><w:jc w:val=`"center`"/></w:pPr><w:r><w:t>Fixed</w:t></w:r></w:p>"
$t.Cell(3, 5).Range.InsertXML($plainFixedXml) | Out-Null

# --- 2. Append a new row for Bug 5 ---
$t.Rows.Add() | Out-Null
$newRow = $t.Rows.Count

$t.Cell($newRow, 1).Range.Text = "5."

$descXml = "<w:p $wNs>" +
  "<w:pPr><w:jc w:val=`"center`"/></w:pPr>" +
  "<w:r><w:t xml:space=`"preserve`">When the shop is selected in the </w:t></w:r>" +
  "<w:proofErr w:type=`"spellStart`"/>" +
  "<w:r><w:t>StartOfDayState</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/>" +
  "<w:r><w:t xml:space=`"preserve`">, it&#8217;s still being selected in </w:t></w:r>" +
  "<w:proofErr w:type=`"spellStart`"/>" +
  "<w:r><w:t>GameplayState</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/>" +
  "</w:p>"
$t.Cell($newRow, 2).Range.InsertXML($descXml) | Out-Null

$t.Cell($newRow, 3).Range.Text = "Yes."
$t.Cell($newRow, 4).Range.Text = "Wei Qi."

# Status cell gets "Fixed" plus the relocated _GoBack bookmark
$statusXml = "<w:p $wNs><w:pPr><w:jc w:val=`"center`"/></w:pPr><w:r><w:t>Fixed</w:t></w:r><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/></w:p>"
$t.Cell($newRow, 5).Range.InsertXML($statusXml) | Out-Null
